# Regenerate save_data to use K (strikeouts) instead of Strike# column.
# This updates the "K" column (column G) values for rows 2-31 on Sheet1,
# replacing the old Strike# counts with the recalculated K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @(
    5,
    1,
    5,
    4,
    1,
    8,
    7,
    2,
    5,
    9,
    7,
    5,
    8,
    7,
    4,
    4,
    9,
    2,
    7,
    6,
    4,
    3,
    8,
    9,
    3,
    1,
    7,
    3,
    6,
    3
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
